# Change drivers for fugitive emissions sectors
# Sheet "Sectors": rows 30-32 (1B1_Fugitive-solid-fuels, 1B2_Fugitive-petr-and-gas,
# 1B2d_Fugitive-other-energy) get a new "activity" driver (column B) and their
# "units" (column C) switch from a literal 1000 number to the text unit "kt".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# Row 31 / 32 share the new driver "refinery-and-natural-gas"; set these before
# row 30's "coal-dom-supply" so the shared-string table gets the same ordering
# as the authored workbook.
$ws.Range("B31").Value = "refinery-and-natural-gas"
$ws.Range("B32").Value = "refinery-and-natural-gas"
$ws.Range("B30").Value = "coal-dom-supply"

# These three cells lose their explicit cell style (revert to the default/Normal
# style) as part of the edit.
$ws.Range("B30").Style = "Normal"
$ws.Range("B31").Style = "Normal"
$ws.Range("B32").Style = "Normal"

# Column C switches from the numeric literal 1000 to the text unit "kt" while
# keeping its existing style.
$ws.Range("C30").Value = "kt"
$ws.Range("C31").Value = "kt"
$ws.Range("C32").Value = "kt"

# Restore the selection state recorded in the saved workbook.
$ws.Range("B30").Select()
